$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, [string]$text)
    $origStyle = $cellRange.Style
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '30.494.33'
Set-TextValue $ws.Range('E2') '  +0.63%  '
Set-TextValue $ws.Range('D3') '2.106.50'
Set-TextValue $ws.Range('E3') '  +1.06%  '
Set-TextValue $ws.Range('D4') '1.007'
Set-TextValue $ws.Range('E4') '  +0.71%  '
Set-TextValue $ws.Range('D5') '334.73'
Set-TextValue $ws.Range('E5') '  +1.87%  '
Set-TextValue $ws.Range('E6') '  +0.75%  '
Set-TextValue $ws.Range('D7') '0.5225'
Set-TextValue $ws.Range('E7') '  -0.03%  '
Set-TextValue $ws.Range('D8') '0.4522'
Set-TextValue $ws.Range('E8') '  +4.68%  '
Set-TextValue $ws.Range('D9') '53.82'
Set-TextValue $ws.Range('E9') '  +14.90%  '
Set-TextValue $ws.Range('D10') '0.08927'
Set-TextValue $ws.Range('E10') '  +1.06%  '
Set-TextValue $ws.Range('D11') '1.182'
Set-TextValue $ws.Range('E11') '  +1.63%  '
Set-TextValue $ws.Range('D12') '24.15'
Set-TextValue $ws.Range('E12') '  -1.43%  '
Set-TextValue $ws.Range('D13') '2.111.56'
Set-TextValue $ws.Range('E13') '  +1.40%  '
Set-TextValue $ws.Range('D14') '6.813'
Set-TextValue $ws.Range('E14') '  +1.10%  '
Set-TextValue $ws.Range('D15') '8.016'
Set-TextValue $ws.Range('E15') '  +3.84%  '
Set-TextValue $ws.Range('D16') '96.74'
Set-TextValue $ws.Range('E16') '  +1.11%  '
Set-TextValue $ws.Range('D17') '1.009'
Set-TextValue $ws.Range('E17') '  +0.81%  '
Set-TextValue $ws.Range('D18') '0.00001141'
Set-TextValue $ws.Range('E18') '  +1.16%  '
Set-TextValue $ws.Range('D19') '0.06666'
Set-TextValue $ws.Range('E19') '  +0.47%  '
Set-TextValue $ws.Range('D20') '19.20'
Set-TextValue $ws.Range('E20') '  +1.56%  '
Set-TextValue $ws.Range('D22') '6.318'
Set-TextValue $ws.Range('E22') '  -0.09%  '
Set-TextValue $ws.Range('D23') '30.556.26'
Set-TextValue $ws.Range('E23') '  +0.62%  '
Set-TextValue $ws.Range('D24') '12.46'
Set-TextValue $ws.Range('E24') '  +0.77%  '
Set-TextValue $ws.Range('D25') '2.353'
Set-TextValue $ws.Range('E25') '  +1.96%  '
Set-TextValue $ws.Range('D26') '2.354.55'
Set-TextValue $ws.Range('E26') '  +1.18%  '
Set-TextValue $ws.Range('D27') '22.25'
Set-TextValue $ws.Range('E27') '  -0.69%  '
Set-TextValue $ws.Range('D28') '162.74'
Set-TextValue $ws.Range('E28') '  +0.47%  '
Set-TextValue $ws.Range('D29') '2.521'
Set-TextValue $ws.Range('E29') '  -2.69%  '
Set-TextValue $ws.Range('D30') '134.11'
Set-TextValue $ws.Range('E30') '  +1.74%  '
Set-TextValue $ws.Range('D31') '1.207'
Set-TextValue $ws.Range('E31') '  +0.52%  '
Set-TextValue $ws.Range('D32') '0.1072'
Set-TextValue $ws.Range('E32') '  +0.09%  '
Set-TextValue $ws.Range('D33') '6.414'
Set-TextValue $ws.Range('E33') '  +3.67%  '
Set-TextValue $ws.Range('E34') '  -1.97%  '
Set-TextValue $ws.Range('E35') '  +1.90%  '
Set-TextValue $ws.Range('D36') '10.42'
Set-TextValue $ws.Range('E36') '  +4.76%  '
Set-TextValue $ws.Range('D37') '5.796'
Set-TextValue $ws.Range('E37') '  +6.27%  '
Set-TextValue $ws.Range('D38') '0.02580'
Set-TextValue $ws.Range('E38') '  +0.38%  '
Set-TextValue $ws.Range('D39') '0.06849'
Set-TextValue $ws.Range('E39') '  +2.51%  '
Set-TextValue $ws.Range('D40') '0.2292'
Set-TextValue $ws.Range('E40') '  +1.14%  '
Set-TextValue $ws.Range('D41') '12.71'
Set-TextValue $ws.Range('E41') '  +0.00%  '
Set-TextValue $ws.Range('D42') '0.6865'
Set-TextValue $ws.Range('E42') '  +0.47%  '
Set-TextValue $ws.Range('E43') '  +0.32%  '
Set-TextValue $ws.Range('B44') 'NEARProtocol'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D44') '2.315'
Set-TextValue $ws.Range('E44') '  +4.86%  '
Set-TextValue $ws.Range('B45') 'EnergySwap'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D45') '14.05'
Set-TextValue $ws.Range('E45') '  +0.21%  '
Set-TextValue $ws.Range('D46') '0.6361'
Set-TextValue $ws.Range('E46') '  -0.46%  '
Set-TextValue $ws.Range('D47') '3.669'
Set-TextValue $ws.Range('E47') '  +1.70%  '
Set-TextValue $ws.Range('B48') 'EOS'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue $ws.Range('D48') '1.250'
Set-TextValue $ws.Range('E48') '  -0.13%  '
Set-TextValue $ws.Range('B49') 'BabyDogeCoin'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D49') '0.00000000348'
Set-TextValue $ws.Range('E49') '  +20.60%  '
Set-TextValue $ws.Range('D50') '1.208'
Set-TextValue $ws.Range('E50') '  +1.37%  '
Set-TextValue $ws.Range('D51') '83.18'
Set-TextValue $ws.Range('E51') '  +1.58%  '
